$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.340.41'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.28%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.846.32'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.21%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9985'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.10'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.30%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6263'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.71%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9984'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.21%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07608'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.31%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2901'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.47%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.74'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.83%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07728'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.28%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.023'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.15%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6780'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.34%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.00001058'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.62%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '82.92'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.99%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.136'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.31%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '29.368.13'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.28%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '227.41'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.93%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.34'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.10%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9983'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.22%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.460'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.09%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9985'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.23%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '158.76'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.86%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.1380'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.63%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.425'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.52%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '17.65'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.21%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.415'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +7.91%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.460'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.70%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.05602'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.63%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.100'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.40%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.067'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.40%  '

$ws.Range("B32").Value = 'ARBITRUM'
$ws.Range("C32").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.162'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.15%  '

$ws.Range("B33").Value = 'LidoDAOToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.832'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.13%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.6965'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.77%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.584'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.15%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.01799'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.09%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.226.50'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.22%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.711'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.56%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.358'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.74%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.8993'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.58%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9982'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.24%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '101.35'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.47%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '65.46'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.28%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.181'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.36%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3986'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.85%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.031'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.05%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.688'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.03%  '

$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.00000000115'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.57%  '

$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.1143'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.53%  '

$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05696'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.30%  '

$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4618'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.20%  '
